$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 4 for the new "CCAR2" coin, shifting the
# existing STEP HERO / WIDI / WSO / Dólar rows down by one.
$ws.Rows.Item(4).Insert()

# Row 2 - BCOIN (values updated, now stored as text with comma decimals)
$ws.Range("B2").Value = "0,89"
$ws.Range("C2").Value = "4,76"

# Row 3 - CCAR (values updated)
$ws.Range("B3").Value = "0,04"
$ws.Range("C3").Value = "0,22"

# Row 4 - CCAR2 (new row)
$ws.Range("A4").Value = "CCAR2"
$ws.Range("B4").Value = "0,0001"
$ws.Range("C4").Value = "0,0005"

# Row 5 - STEP HERO (values updated)
$ws.Range("B5").Value = "0,09"
$ws.Range("C5").Value = "0,48"

# Row 6 - WIDI (values updated)
$ws.Range("B6").Value = "0,10"
$ws.Range("C6").Value = "0,54"

# Row 7 - WSO (values updated)
$ws.Range("B7").Value = "0,01"
$ws.Range("C7").Value = "0,07"

# Row 8 - Dólar (value updated; C8 keeps the original empty text cell
# that shifted down from C7 via the row insert above, so it is left
# untouched here)
$ws.Range("B8").Value = "5,26"
